# Update column F ("dSF") values for the specified rows as per the
# repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 9
    3  = 2
    4  = 1
    5  = 8
    7  = -2
    8  = 2
    9  = -1
    10 = -6
    11 = 5
    12 = -2
    13 = 4
    14 = -1
    15 = 7
    16 = 2
    17 = -3
    18 = -1
    19 = -1
    20 = -3
    21 = 2
    23 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
